# chart updated with latest stats
# Appends the 2020-04-14 daily COVID-19 stats block (rows 23-29) to Sheet1,
# mirroring the existing per-date layout:
#   Col A = City name, Col B = City-wise count, Col C = Date,
#   Col D = Overall-count header label, Col E = Overall count value

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $val) {
    # Force the cell to be stored as text (matching every other data cell in
    # this sheet, which are all shared-string / General-style cells) instead
    # of letting Excel auto-convert numeric- or date-looking strings into
    # numbers/dates with an applied number format.
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

$newRows = @(
    @{ Row = 23; A = "ICT";         B = "131";   C = "2020-04-14"; D = "Recovered";       E = "1,378" },
    @{ Row = 24; A = "Punjab";      B = "2,881"; C = "2020-04-14"; D = "Critical";         E = "46" },
    @{ Row = 25; A = "Sindh";       B = "1,518"; C = "2020-04-14"; D = "Deaths";           E = "96" },
    @{ Row = 26; A = "KP";          B = "800";   C = "2020-04-14"; D = "Cases (24 HRS)";   E = "121" },
    @{ Row = 27; A = "Balochistan"; B = "231";   C = "2020-04-14"; D = "Deaths (24 HRS)";  E = "3" },
    @{ Row = 28; A = "AJK";         B = "43";    C = "2020-04-14"; D = "Tests (24 HRS)";   E = "3,157" },
    @{ Row = 29; A = "GB";          B = "233";   C = "2020-04-14"; D = "Total Tests";      E = "69,928" }
)

# Write column-by-column (all of A, then all of B, ...) rather than
# row-by-row so that new entries land in the shared-strings table in the
# same order the original authoring tool produced them in.
foreach ($r in $newRows) { Set-TextCell $ws.Range("A$($r.Row)") $r.A }
foreach ($r in $newRows) { Set-TextCell $ws.Range("B$($r.Row)") $r.B }
foreach ($r in $newRows) { Set-TextCell $ws.Range("C$($r.Row)") $r.C }
foreach ($r in $newRows) { Set-TextCell $ws.Range("D$($r.Row)") $r.D }
foreach ($r in $newRows) { Set-TextCell $ws.Range("E$($r.Row)") $r.E }
